# Weekly data refresh: insert 3 new "Black Amber" price rows (date 2022-01-04)
# right after the "Angeleno" block ending at row 27, shifting every
# subsequent record down by 3 rows (old row 28 -> new row 31, ...,
# old row 116 -> new row 119). Sheet dimension grows from A1:T116 to A1:T119.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 28..116 down by 3 rows, preserving their content/formatting.
$ws.Rows("28:30").Insert()

# New records to populate into the freshly opened rows 28-30.
$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44565, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Especial", 56, 14000, 14000, 14000, "`$/caja 15 kilos empedrada", "Provincia de San Felipe de Aconcagua", 933, 15),
    @(3, "Femacal de La Calera", "Coquimbo", 44565, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 50, 12000, 12000, 12000, "`$/caja 15 kilos empedrada", "Provincia de San Felipe de Aconcagua", 800, 15),
    @(3, "Femacal de La Calera", "Coquimbo", 44565, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Segunda", 48, 10000, 10000, 10000, "`$/caja 15 kilos empedrada", "Provincia de San Felipe de Aconcagua", 667, 15)
)

$startRow = 28
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowVals = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
